$wb = $excel.ActiveWorkbook

# Add the new "Address Book" sheet right after the "Review" sheet
$reviewSheet = $wb.Worksheets.Item("Review")
$ws = $wb.Worksheets.Add($null, $reviewSheet)
$ws.Name = "Address Book"

# Header row (row 1)
$ws.Range("A1").Value = "DataSet"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Prod UserName"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Password"
$ws.Range("F1").Value = "Confirm Password"
$ws.Range("G1").Value = "FirstName"
$ws.Range("H1").Value = "LastName"
$ws.Range("I1").Value = "Street"
$ws.Range("J1").Value = "City"
$ws.Range("K1").Value = "Region"
$ws.Range("L1").Value = "postcode"
$ws.Range("M1").Value = "phone"
$ws.Range("N1").Value = "Products"
$ws.Range("O1").Value = "Color"
$ws.Range("P1").Value = "Size"
$ws.Range("Q1").Value = "Colorproduct"
$ws.Range("R1").Value = "Quantity"
$ws.Range("S1").Value = "methods"
$ws.Range("T1").Value = "cardNumber"
$ws.Range("U1").Value = "ExpMonthYear"
$ws.Range("V1").Value = "cvv"
$ws.Range("A1:V1").Interior.Color = 65535

# Row 2
$ws.Range("A2").Value = "Account"
$ws.Range("B2").Value = "testersemail.278@gmail.com"
$ws.Range("D2").Value = "testersemail.278@gmail.com"
$ws.Range("E2").Value = "Testers@278"
$ws.Range("F2").Value = "Testers@278"
$ws.Range("G2").Value = "QA"
$ws.Range("H2").Value = "TEST"
$ws.Range("I2").Value = "6 Walnut Valley Dr"
$ws.Range("J2").Value = "Little Rock"
$ws.Range("K2").Value = "Arkansas"
$ws.Range("L2").Value = 72211
$ws.Range("M2").Value = 9898989898

# Row 3
$ws.Range("A3").Value = "AddressBook"
$ws.Range("G3").Value = "QA"
$ws.Range("H3").Value = "TEST"
$ws.Range("I3").Value = "844 N Colony Rd"
$ws.Range("J3").Value = "Wallingford"
$ws.Range("K3").Value = "Connecticut"
$ws.Range("L3").Value = "'06492"
$ws.Range("M3").Value = 9898989898

# Autofit the columns like Excel does after data entry
$ws.Columns("A:V").AutoFit()

# Select A4 (below the data) to match final cursor position, then make this the active/visible tab
$ws.Range("A4").Select()

# Clear the stale selection on the Review sheet and select the default cell
$reviewSheet.Range("A1").Select()

# Put focus back on the new Address Book sheet, as the workbook's active sheet
$ws.Activate()

# Update selection on "Checkout payments": select rows 1:2 (whole rows)
$cp = $wb.Worksheets.Item("Checkout payments")
$cp.Rows("1:2").Select()

# Leave Address Book as the final active sheet/tab
$ws.Activate()
